$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (row => B,C,D,E,F,G)
$data = @{
    2 = @(0.6545652718822623, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 0, 1.378593816457939)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    4 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    5 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 3.536033448013082)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G
}
